# adapt slides for CS 536 in the spring
#
# 1. The cached "today" text baked into the auto-updating date field
#    (type="datetimeFigureOut") on the slide master and every slide
#    layout gets refreshed from 11/18/2024 -> 12/20/2024 (this is what
#    PowerPoint does to every such field whenever the deck is saved).
# 2. Slide 1 ("Note to CS 4536 Students") is marked as a hidden slide
#    (show="0") so it is skipped during the slide show.

$p = $ppt.ActivePresentation

$oldDate = "11/18/2024"
$newDate = "12/20/2024"

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = -1
        try {
            $phType = $shp.PlaceholderFormat.Type
        } catch {
            $phType = -1
        }
        if ($phType -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster

# Slide master's own date placeholder.
Update-DatePlaceholder $master.Shapes

# Every slide layout has its own cached copy of the date placeholder.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Hide the first slide ("Note to CS 4536 Students") from the slide show.
$s1 = $p.Slides.Item(1)
$s1.SlideShowTransition.Hidden = $true
